# Adds 8 new printable parameters (fan/horn/adapter-wall dimensions) to
# the "Sheet1" parameter table, tweaks two existing values, and leaves the
# selection on the updated B27 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# armslotb (B23): 6 -> 4
$ws.Range("B23").Value = 4

# distancefromleft (B27): =B4 - 18 -> =B4 - 20
$ws.Range("B27").Formula = "=B4 - 20"

# Append 8 new parameter rows (37-44) below the existing table.
$newParams = @(
    @{ Name = "fanhb";             Value = 11 },
    @{ Name = "fanhs";             Value = 9.6 },
    @{ Name = "adapterwallthick";  Value = 6 },
    @{ Name = "adapterwallslotd";  Value = 2 },
    @{ Name = "armhornthick";      Value = 3 },
    @{ Name = "hornthick";         Value = 1.5 },
    @{ Name = "hornscrewholed";    Value = 2 },
    @{ Name = "hornscrewheadd";    Value = 3 }
)

$row = 37
foreach ($param in $newParams) {
    $ws.Cells.Item($row, 1).Value = $param.Name
    $ws.Cells.Item($row, 2).Value = $param.Value
    $row++
}

# Match the saved selection state (was A37, now B27).
$ws.Range("B27").Select()
